$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 86) with the next business day's series data.
#
# Column A holds a date written as plain text (e.g. "04-11-2021" in the row
# above). Simply assigning the literal string "05-11-2021" to a cell's
# .Value makes Excel's smart data-entry parser recognize it as a valid
# MM-DD-YYYY date and silently convert the cell to a date serial number
# (adding a date number-format style in the process). To keep the value as
# genuine text -- matching how the rest of the column is stored -- build it
# as a formula result (which Excel does not reinterpret as a date) and then
# flatten that formula down to a static value with a values-only paste.
$ws.Range("A86").Formula = "=SUBSTITUTE(""13-11-2021"",""13"",""05"")"
$ws.Range("A86").Copy()
$ws.Range("A86").PasteSpecial(-4163)

$ws.Range("B86").Value = 50000
$ws.Range("C86").Value = 170000
$ws.Range("D86").Value = 50000
$ws.Range("E86").Value = 15000
$ws.Range("F86").Value = 35000
$ws.Range("G86").Value = 3.19
